$d = $word.ActiveDocument

# Locate the last paragraph (the one ending with "The key to this problem...")
$lastParaIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($lastParaIndex)

# Move the insertion point to the very end of that paragraph's text (before the
# paragraph mark) and insert a paragraph break plus the new content there.
$insertRange = $lastPara.Range
$insertRange.Collapse(0)  # wdCollapseEnd
# Move back one character to land before the paragraph mark
$insertRange.MoveEnd(1, -1) | Out-Null
$insertRange.Collapse(0)

$insertRange.InsertParagraphAfter()
$insertRange.Collapse(0)
$insertRange.InsertParagraphAfter()
$insertRange.Collapse(0)
$insertRange.InsertParagraphAfter()
$insertRange.Collapse(0)

# Now find paragraph markers to fill in text.
$total = $d.Paragraphs.Count
Write-Host "Total paragraphs after insert:" $total

$pEmpty = $d.Paragraphs.Item($total - 2)
$pHeading = $d.Paragraphs.Item($total - 1)
$pBody = $d.Paragraphs.Item($total)

$pHeading.Range.Text = "2) Break the Problem Apart"
$pBody.Range.Text = "The problem with this problem is that the girl does not use a method of counting that utilizes each finger equally. If she counted 1-5 from thumb to pinky and 6-10 from pinky to thumb (or thumb to pinky again), we could easily divide each number by 5 and figure out which finger the number would land on. Instead we have to figure out what equation will work in the same manner with her odd counting style. If we can figure this out on a small scale, we should then be able to apply it on a larger scale."
